$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.074.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.452.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.83%  '

$ws.Range("E7").Value = '  +0.94%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.30%  '

$ws.Range("E10").Value = '  +1.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0806'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.122'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.24'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.836.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.434.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.840'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.931.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0935'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.99%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '246.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.23%  '

$ws.Range("E24").Value = '  +2.29%  '

$ws.Range("E25").Value = '  +0.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.08%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.130'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.34'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.49%  '

$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0762'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.25%  '

$ws.Range("E38").Value = '  -0.32%  '

$ws.Range("E39").Value = '  +0.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '124.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.20%  '

$ws.Range("E41").Value = '  +2.24%  '

$ws.Range("E42").Value = '  +1.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0293'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.94%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.962.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.72%  '

$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("E47").Value = '  -2.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '77.67'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.89%  '
